# Update "想去人数" (want-to-go count) figures in the F column, reflecting a
# fresh scrape of the source site. Two sheets list overlapping events and
# both need their F-column counters bumped to the newer numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 156
$ws1.Range("F4").Value = 73
$ws1.Range("F5").Value = 517
$ws1.Range("F6").Value = 1550
$ws1.Range("F7").Value = 3
$ws1.Range("F8").Value = 1180
$ws1.Range("F9").Value = 115
$ws1.Range("F10").Value = 218
$ws1.Range("F11").Value = 157
$ws1.Range("F12").Value = 2
$ws1.Range("F14").Value = 3
$ws1.Range("F15").Value = 229
$ws1.Range("F16").Value = 127
$ws1.Range("F17").Value = 189
$ws1.Range("F18").Value = 176

# --- Sheet "全部类型" (all event types) --------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 156
$ws4.Range("F4").Value = 73
$ws4.Range("F5").Value = 517
$ws4.Range("F6").Value = 1551
$ws4.Range("F8").Value = 3
$ws4.Range("F9").Value = 1180
$ws4.Range("F10").Value = 115
$ws4.Range("F11").Value = 218
$ws4.Range("F12").Value = 157
$ws4.Range("F13").Value = 2
$ws4.Range("F15").Value = 3
$ws4.Range("F16").Value = 229
$ws4.Range("F17").Value = 127
$ws4.Range("F18").Value = 189
$ws4.Range("F19").Value = 176
